$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.023.16'
$ws.Range('E2').Value = '  -1.77%  '
$ws.Range('D3').Value = '2.971.31'
$ws.Range('E3').Value = '  -0.36%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '595.81'
$ws.Range('E5').Value = '  +2.94%  '
$ws.Range('D6').Value = '141.95'
$ws.Range('E6').Value = '  -2.20%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  -0.89%  '
$ws.Range('D9').Value = '2.971.79'
$ws.Range('E9').Value = '  -0.37%  '
$ws.Range('E10').Value = '  -1.60%  '
$ws.Range('E11').Value = '  +6.30%  '
$ws.Range('E12').Value = '  +3.11%  '
$ws.Range('E13').Value = '  -0.08%  '
$ws.Range('D14').Value = '33.96'
$ws.Range('E14').Value = '  -1.27%  '
$ws.Range('E15').Value = '  +2.32%  '
$ws.Range('D16').Value = '3.462.98'
$ws.Range('E16').Value = '  -0.29%  '
$ws.Range('E17').Value = '  -2.10%  '
$ws.Range('D18').Value = '61.011.37'
$ws.Range('E18').Value = '  -1.70%  '
$ws.Range('D19').Value = '2.972.38'
$ws.Range('E19').Value = '  -0.47%  '
$ws.Range('D20').Value = '446.42'
$ws.Range('E20').Value = '  -1.64%  '
$ws.Range('D21').Value = '14.08'
$ws.Range('E21').Value = '  +2.11%  '
$ws.Range('D22').Value = '0.677'
$ws.Range('E22').Value = '  +0.53%  '
$ws.Range('E23').Value = '  +0.18%  '
$ws.Range('D24').Value = '81.86'
$ws.Range('E24').Value = '  +2.65%  '
$ws.Range('E25').Value = '  -4.74%  '
$ws.Range('D26').Value = '10.30'
$ws.Range('E26').Value = '  +3.60%  '
$ws.Range('D27').Value = '11.84'
$ws.Range('E27').Value = '  -2.39%  '
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('E29').Value = '  +2.91%  '
$ws.Range('E30').Value = '  +0.02%  '
$ws.Range('E31').Value = '  -1.64%  '
$ws.Range('E32').Value = '  -1.84%  '
$ws.Range('D33').Value = '27.00'
$ws.Range('E33').Value = '  +1.31%  '
$ws.Range('D35').Value = '0.0₃0807'
$ws.Range('E35').Value = '  +3.75%  '
$ws.Range('E36').Value = '  -0.79%  '
$ws.Range('E37').Value = '  +0.45%  '
$ws.Range('D38').Value = '50.19'
$ws.Range('E38').Value = '  +0.52%  '
$ws.Range('D39').Value = '2.04'
$ws.Range('E39').Value = '  -2.53%  '
$ws.Range('D40').Value = '8.98'
$ws.Range('E40').Value = '  +0.28%  '
$ws.Range('E41').Value = '  +9.48%  '
$ws.Range('E42').Value = '  -2.17%  '
$ws.Range('D43').Value = '389.86'
$ws.Range('E43').Value = '  -3.89%  '
$ws.Range('D44').Value = '39.00'
$ws.Range('E44').Value = '  +1.82%  '
$ws.Range('E45').Value = '  -0.45%  '
$ws.Range('E46').Value = '  -3.81%  '
$ws.Range('D47').Value = '2.676.60'
$ws.Range('E47').Value = '  -2.86%  '
$ws.Range('D48').Value = '129.92'
$ws.Range('E48').Value = '  +2.20%  '
$ws.Range('E49').Value = '  +0.08%  '
$ws.Range('E50').Value = '  -0.57%  '
$ws.Range('E51').Value = '  -0.39%  '
